$d = $word.ActiveDocument

$chomageUrl = "https://fr.countryeconomy.com/marche-du-travail/chomage"
$visaUrl    = "https://visaguide.world/visa-free-countries/french-passport/"

# --- Stage 1: split the single "chomage" paragraph into two paragraphs ---
# using unique placeholder markers, via Find/Replace with a paragraph mark
# (^p) in the replacement. This produces clean paragraph splits without
# leaving stray empty runs behind.
$d.Content.Find.Execute($chomageUrl, $false, $false, $false, $false, $false, `
    $true, 1, $false, "CHOMAGE_MARKER^pVISA_MARKER", 2) | Out-Null

# --- Stage 2: insert a clean empty paragraph right after the (future)
# hyperlink paragraph, and another clean empty paragraph right after the
# (future) visa paragraph (this becomes the final trailing blank line). ---
$d.Content.Find.Execute("CHOMAGE_MARKER", $false, $false, $false, $false, `
    $false, $true, 1, $false, "CHOMAGE_MARKER^p", 2) | Out-Null
$d.Content.Find.Execute("VISA_MARKER", $false, $false, $false, $false, `
    $false, $true, 1, $false, "VISA_MARKER^p", 2) | Out-Null

# At this point the paragraph list (tail) looks like:
#  ... numbeo hyperlink / empty / CHOMAGE_MARKER / empty / VISA_MARKER / empty

# --- Stage 3: turn the CHOMAGE_MARKER paragraph into the real hyperlink ---
$chomagePara = $d.Paragraphs(13)
$chomageRange = $chomagePara.Range
$chomageRange.MoveEnd(1, -1) | Out-Null
$chomageRange.Text = $chomageUrl

$chomagePara = $d.Paragraphs(13)
$chomageRange = $chomagePara.Range
$chomageRange.MoveEnd(1, -1) | Out-Null
$link = $d.Hyperlinks.Add($chomageRange, $chomageUrl, $null, $null, $chomageUrl)
$link.Range.Style = "Lienhypertexte"

# --- Stage 4: turn the VISA_MARKER paragraph into the visaguide line,
# made of two separate runs: the URL, and "    lien visa" ---
$visaPara = $d.Paragraphs(15)
$visaRange = $visaPara.Range
$visaRange.MoveEnd(1, -1) | Out-Null
$visaRange.Text = $visaUrl

$visaPara = $d.Paragraphs(15)
$visaRange = $visaPara.Range
$visaRange.MoveEnd(1, -1) | Out-Null
$visaRange.Collapse(0) | Out-Null
$visaRange.InsertAfter("    lien visa") | Out-Null

# force the appended text to live in its own run (rather than being
# silently merged back into the preceding run)
$secondRunStart = $visaRange.Start + $visaUrl.Length
$secondRunEnd = $visaRange.Start + $visaUrl.Length + "    lien visa".Length
$secondRun = $d.Range($secondRunStart, $secondRunEnd)
$secondRun.Bold = 1
$secondRun.Bold = 0

"done"
